$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from the generic pattern name to the
# descriptive "Venta - Plan" title used once the sheet is finalized
# with this period's actual sale figures.
$ws.Name = "Venta - Plan"
